# Solicitar tutoría - estudiante
# Collapse the four separate "paralelo_a/b/c/d" flag columns (H:K) into a
# single "paralelo" column holding a comma separated list of parallels
# (e.g. "A,B"), and move the trailing "ciclo" column (L) left into I so the
# used range shrinks from A1:L13 down to A1:I13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the now-unused J:L columns (old paralelo_b/paralelo_c/paralelo_d
# + the old ciclo column) so the sheet's used range collapses to column I.
$ws.Range("J1:L13").ClearContents()

# Rename the remaining headers: H becomes the combined "paralelo" column,
# and I becomes "ciclo" (previously the header text that lived in L1).
$ws.Range("H1").Value = "paralelo"
$ws.Range("I1").Value = "ciclo"

# Every docente row now reports the parallels they teach as "A,B" in the
# paralelo column, keeping "NA" in the ciclo column.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value = "A,B"
    $ws.Cells.Item($r, 9).Value = "NA"
}

# Match the author's final cursor position recorded in the sheet view.
$ws.Range("G18").Select()
